# TM_mujeres.xlsx — "agregando ejercicio de pensiones y longevidad"
#
# 1) Fix seven B-column figures that were mis-keyed with a stray decimal
#    point (e.g. 994.4 instead of 994400) so they read as whole pesos,
#    consistent with the rest of the l(x) column.
# 2) Clear the leftover direct number-formatting on the C:E columns
#    (d(x), q(x), e0(x)) — they only need the workbook's General format.
# 3) Leave the selection on column F, as the author did before saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- 1) Correct the mis-entered l(x) values in column B ------------------
$ws.Range("B18").Value = 994400
$ws.Range("B33").Value = 982452
$ws.Range("B64").Value = 792450
$ws.Range("B77").Value = 411435
$ws.Range("B81").Value = 253452
$ws.Range("B82").Value = 215411
$ws.Range("B93").Value = 4425

# --- 2) Strip the direct number format from C2:E97 -----------------------
$ws.Range("C2:E97").ClearFormats()

# --- 3) Match the saved selection (whole column F) ------------------------
$ws.Columns("F").Select() | Out-Null
